# Applies updated FFXIV market/leve profit data pulled by the scheduled runner.
# Generated from the authoritative diff of the workbook's OOXML.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1674.1052
$ws.Range("I98").Value = 1794.375
$ws.Range("J98").Value = 1032.6666
$ws.Range("K98").Value = 1794.375
$ws.Range("L98").Value = 1032.6666
$ws.Range("M98").Value = -296.375
$ws.Range("N98").Value = -4028.6666

$ws.Range("H116").Value = 4114.524
$ws.Range("I116").Value = 4071.3572
$ws.Range("J116").Value = 4200.857
$ws.Range("K116").Value = 4071.3572
$ws.Range("L116").Value = 4200.857
$ws.Range("M116").Value = -629.3571999999999
$ws.Range("N116").Value = -11084.857

$ws.Range("H122").Value = 1674.1052
$ws.Range("I122").Value = 1794.375
$ws.Range("J122").Value = 1032.6666
$ws.Range("K122").Value = 5383.125
$ws.Range("L122").Value = 3097.9998
$ws.Range("M122").Value = -2933.125
$ws.Range("N122").Value = -7997.9998

$ws.Range("H132").Value = 1891.0244
$ws.Range("I132").Value = 1925.6389
$ws.Range("J132").Value = 1641.8
$ws.Range("K132").Value = 5776.9167
$ws.Range("L132").Value = 4925.4
$ws.Range("M132").Value = -3246.9167
$ws.Range("N132").Value = -9985.4

$ws.Range("H137").Value = 908.13336
$ws.Range("I137").Value = 840.1539
$ws.Range("J137").Value = 1350
$ws.Range("K137").Value = 2520.4617
$ws.Range("L137").Value = 4050
$ws.Range("M137").Value = 29.53830000000016
$ws.Range("N137").Value = -9150

$ws.Range("H138").Value = 1107.09
$ws.Range("I138").Value = 619.2461499999999
$ws.Range("J138").Value = 2013.0857
$ws.Range("K138").Value = 1857.73845
$ws.Range("L138").Value = 6039.257100000001
$ws.Range("M138").Value = 3282.26155
$ws.Range("N138").Value = -16319.2571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3524445.8
$ws.Range("I32").Value = 4403733
$ws.Range("J32").Value = 7297.5264
$ws.Range("K32").Value = 4403733
$ws.Range("L32").Value = 7297.5264
$ws.Range("M32").Value = -4403446
$ws.Range("N32").Value = -7871.5264

$ws.Range("H45").Value = 1402.579
$ws.Range("I45").Value = 1197.9333
$ws.Range("J45").Value = 2170
$ws.Range("K45").Value = 1197.9333
$ws.Range("L45").Value = 2170
$ws.Range("M45").Value = -820.9332999999999
$ws.Range("N45").Value = -2924

$ws.Range("H110").Value = 576.5714
$ws.Range("I110").Value = 532.4
$ws.Range("J110").Value = 687
$ws.Range("K110").Value = 532.4
$ws.Range("L110").Value = 687
$ws.Range("M110").Value = 1512.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3940.0435
$ws.Range("I105").Value = 4837.857
$ws.Range("J105").Value = 2543.4443
$ws.Range("K105").Value = 4837.857
$ws.Range("L105").Value = 2543.4443
$ws.Range("M105").Value = -3090.857
$ws.Range("N105").Value = -6037.4443

$ws.Range("H132").Value = 139078
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 139078
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 139078
$ws.Range("N132").Value = -149198

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5168.5
$ws.Range("I16").Value = 4402.2
$ws.Range("J16").Value = 9000
$ws.Range("K16").Value = 4402.2
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = -4115.2
$ws.Range("N16").Value = -9574

$ws.Range("H31").Value = 1485.175
$ws.Range("I31").Value = 1475.1082
$ws.Range("J31").Value = 1609.3334
$ws.Range("K31").Value = 1475.1082
$ws.Range("L31").Value = 1609.3334
$ws.Range("M31").Value = -1180.1082
$ws.Range("N31").Value = -2199.3334

$ws.Range("H34").Value = 1485.175
$ws.Range("I34").Value = 1475.1082
$ws.Range("J34").Value = 1609.3334
$ws.Range("K34").Value = 1475.1082
$ws.Range("L34").Value = 1609.3334
$ws.Range("M34").Value = -1273.1082
$ws.Range("N34").Value = -2013.3334

$ws.Range("H113").Value = 5168.5
$ws.Range("I113").Value = 4402.2
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 4402.2
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -2232.2
$ws.Range("N113").Value = -13340

$ws.Range("H122").Value = 994.36365
$ws.Range("I122").Value = 820.6667
$ws.Range("J122").Value = 1202.8
$ws.Range("K122").Value = 2462.0001
$ws.Range("L122").Value = 3608.4
$ws.Range("M122").Value = -12.0001000000002
$ws.Range("N122").Value = -8508.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1035.8182
$ws.Range("I5").Value = 921.55554
$ws.Range("J5").Value = 1550
$ws.Range("K5").Value = 2764.66662
$ws.Range("L5").Value = 4650
$ws.Range("M5").Value = -2652.66662
$ws.Range("N5").Value = -4874

$ws.Range("H109").Value = 1533.3334
$ws.Range("I109").Value = 800
$ws.Range("J109").Value = 3000
$ws.Range("K109").Value = 2400
$ws.Range("L109").Value = 9000
$ws.Range("M109").Value = -1360

$ws.Range("H122").Value = 294592.16
$ws.Range("I122").Value = 275.66666
$ws.Range("J122").Value = 357659.97
$ws.Range("K122").Value = 2480.99994
$ws.Range("L122").Value = 3218939.73
$ws.Range("M122").Value = -30.9999399999997
$ws.Range("N122").Value = -3223839.73

$ws.Range("H131").Value = 764.41
$ws.Range("I131").Value = 347.6875
$ws.Range("J131").Value = 843.7857
$ws.Range("K131").Value = 1043.0625
$ws.Range("L131").Value = 2531.3571
$ws.Range("M131").Value = 3996.9375
$ws.Range("N131").Value = -12611.3571

$ws.Range("H135").Value = 1035.8182
$ws.Range("I135").Value = 921.55554
$ws.Range("J135").Value = 1550
$ws.Range("K135").Value = 8293.99986
$ws.Range("L135").Value = 13950
$ws.Range("M135").Value = -5758.99986
$ws.Range("N135").Value = -19020

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2077.4236
$ws.Range("I132").Value = 1736.6938
$ws.Range("J132").Value = 2541.1943
$ws.Range("K132").Value = 5210.0814
$ws.Range("L132").Value = 7623.5829
$ws.Range("M132").Value = -2680.0814
$ws.Range("N132").Value = -12683.5829

$ws.Range("H136").Value = 3057.4707
$ws.Range("I136").Value = 1489.2273
$ws.Range("J136").Value = 12915
$ws.Range("K136").Value = 4467.6819
$ws.Range("L136").Value = 38745
$ws.Range("M136").Value = -1917.6819
$ws.Range("N136").Value = -43845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 25000
$ws.Range("I76").Value = 25000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 25000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -24685

$ws.Range("H79").Value = 25000
$ws.Range("I79").Value = 25000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 25000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -23908

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
